$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L6").Value = 506.88
$ws1.Range("L24").Value = "1 de 22"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F6").Value = 1265.01
$ws2.Range("F24").Value = 18923.72

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D11").Value = 506.88
$ws3.Range("E11").Value = 2415.34458185274
$ws3.Range("F11").Value = 0.1734568941578848

$ws3.Range("D15").Value = 18923.72
$ws3.Range("E15").Value = 39279.74623249458
$ws3.Range("F15").Value = 0.325130464299307
